# Update the date heading
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-04-18 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-04-19 Friday", 2)

# Update the division-problem table cells. Several source strings repeat
# (e.g. "63÷2=31, 1" appears more than once) but map to different targets,
# so each cell must be addressed individually by (row, column) rather than
# via a global text replace.
$tbl = $d.Tables.Item(1)

$values = @(
    @(1,  1, "48÷7=6, 6"),
    @(1,  2, "10÷9=1, 1"),
    @(1,  3, "61÷4=15, 1"),
    @(1,  4, "80÷8=10, 0"),
    @(1,  5, "15÷2=7, 1"),

    @(5,  1, "34÷2=17, 0"),
    @(5,  2, "94÷4=23, 2"),
    @(5,  3, "42÷4=10, 2"),
    @(5,  4, "75÷2=37, 1"),
    @(5,  5, "88÷6=14, 4"),

    @(9,  1, "33÷6=5, 3"),
    @(9,  2, "94÷5=18, 4"),
    @(9,  3, "18÷8=2, 2"),
    @(9,  4, "33÷4=8, 1"),
    @(9,  5, "76÷4=19, 0"),

    @(13, 1, "91÷8=11, 3"),
    @(13, 2, "80÷3=26, 2"),
    @(13, 3, "44÷2=22, 0"),
    @(13, 4, "76÷5=15, 1"),
    @(13, 5, "35÷5=7, 0"),

    @(17, 1, "60÷9=6, 6"),
    @(17, 2, "73÷7=10, 3"),
    @(17, 3, "56÷9=6, 2"),
    @(17, 4, "35÷7=5, 0"),
    @(17, 5, "27÷3=9, 0")
)

foreach ($item in $values) {
    $row = $item[0]
    $col = $item[1]
    $text = $item[2]
    $tbl.Cell($row, $col).Range.Text = $text
}
